$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Paragraph 2: turn the field `{ m:userdoc 'zone1' }` into four plain
# text runs: "{", "m", ":userdoc 'zone1'", "}"
# ------------------------------------------------------------------
$field1 = $d.Fields.Item(1)
$p2Start = $d.Paragraphs.Item(2).Range.Start
$field1.Delete()

$cursor = $d.Range($p2Start, $p2Start)
$cursor.InsertAfter("{")
$cursor.Collapse(0)
$d.Bookmarks.Add("tmpSplitA", $cursor) | Out-Null
$cursor.InsertAfter("m")
$cursor.Collapse(0)
$d.Bookmarks.Add("tmpSplitB", $cursor) | Out-Null
$cursor.InsertAfter(":userdoc 'zone1'")
$cursor.Collapse(0)
$d.Bookmarks.Add("tmpSplitC", $cursor) | Out-Null
$cursor.InsertAfter("}")

$d.Bookmarks("tmpSplitA").Delete()
$d.Bookmarks("tmpSplitB").Delete()
$d.Bookmarks("tmpSplitC").Delete()

# ------------------------------------------------------------------
# Paragraph 4: turn the field `{ m:enduserdoc }` into two plain text
# runs "{m:" and "enduserdoc}" while keeping the _GoBack bookmark
# between them.
# ------------------------------------------------------------------
$field2 = $d.Fields.Item(1)
$p4Start = $d.Paragraphs.Item(4).Range.Start
$field2.Delete()

$cursor2 = $d.Range($p4Start, $p4Start)
$cursor2.InsertAfter("{m:")
$cursor2.Collapse(0)
$d.Bookmarks.Add("tmpSplitD", $cursor2) | Out-Null
$cursor2.InsertAfter("enduserdoc}")
$d.Bookmarks("tmpSplitD").Delete()

$goBackRange = $d.Range($p4Start + 3, $p4Start + 3)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
